# Improve client library API diagram.
#
# Nudges a handful of boxes in the "micro-ROS client library" diagram on
# slide 3 (corner-shaped title band + the small feature boxes anchored to
# it), reshapes the "corner" connector band, and restores the intended
# (brand) theme colors on the slide master so the diagram renders with its
# proper palette instead of the generic Office defaults.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# EMU -> point helper. Shape.Left/.Top are stored as a lower-precision
# (points, ~4 decimal digits) value internally, so a straight EMU/12700
# division can land a single EMU short of the intended target after the
# round-trip. Nudging by half an EMU (in points) keeps the conversion on
# the correct side of that rounding boundary without perturbing the
# intended position in any visible way.
$EmuPerPt = 12700
$Eps = 0.5 / $EmuPerPt
function ToPt([double]$emu) { return ($emu / $EmuPerPt) + $Eps }

# --- Reposition the feature boxes that sit on top of the "corner" band ---

$moves = @(
    @{ Name = "Google Shape;162;g6cd1c02217_0_0"; X = 4446138; Y = 1076225 },
    @{ Name = "Google Shape;166;g6cd1c02217_0_0"; X = 5407775; Y = 1232200 },
    @{ Name = "Google Shape;167;g6cd1c02217_0_0"; X = 6244102; Y = 1232200 },
    @{ Name = "Google Shape;168;g6cd1c02217_0_0"; X = 7080450; Y = 1232200 },
    @{ Name = "Google Shape;169;g6cd1c02217_0_0"; X = 5407775; Y = 1821200 },
    @{ Name = "Google Shape;184;g6cd1c02217_0_0"; X = 6244113; Y = 1821200 }
)

foreach ($m in $moves) {
    $sh = $s.Shapes.Item($m.Name)
    $sh.Left = ToPt($m.X)
    $sh.Top = ToPt($m.Y)
}

# --- Reshape the "corner" connector band behind the title (adj2 handle) ---

$corner = $s.Shapes.Item("Google Shape;161;g6cd1c02217_0_0")
$corner.Adjustments(2) = 2.67408

# --- Restore the brand color theme on the slide master ---
# (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, in that order)

$tcs = $s.ThemeColorScheme
$tcs.Item(1).RGB  = 0         # dk1      #000000
$tcs.Item(2).RGB  = 16777215  # lt1      #FFFFFF
$tcs.Item(3).RGB  = 5800213   # dk2      #158158
$tcs.Item(4).RGB  = 15987699  # lt2      #F3F3F3
$tcs.Item(5).RGB  = 13077765  # accent1  #058DC7
$tcs.Item(6).RGB  = 3322960   # accent2  #50B432
$tcs.Item(7).RGB  = 1791725   # accent3  #ED561B
$tcs.Item(8).RGB  = 61421     # accent4  #EDEF00
$tcs.Item(9).RGB  = 15059748  # accent5  #24CBE5
$tcs.Item(10).RGB = 7529828   # accent6  #64E572
$tcs.Item(11).RGB = 13369378  # hlink    #2200CC
$tcs.Item(12).RGB = 9116245   # folHlink #551A8B
